# "Update Waste model version": reorder the WASTE emission-factor entries
# in the Lists sheet's Emission column (E3:E5) from
#   CH4_WASTE, CO2e_WASTE, N2O_WASTE, CO2_WASTE
# to
#   CH4_WASTE, N2O_WASTE, CO2_WASTE, CO2e_WASTE

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lists")

$ws.Range("E3").Value = "N2O_WASTE"
$ws.Range("E4").Value = "CO2_WASTE"
$ws.Range("E5").Value = "CO2e_WASTE"
